# fall 22 week 11 complete plus 9 ball skill level evals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matchup results for week 11
$ws.Range("C3").Value = "0/3"
$ws.Range("B5").Value = "3/0"
